$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp refresh (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 04:35"

# Row 5: updated totals for this country
$ws.Range("B5").Value = 349113
$ws.Range("C5").Value = 1715
$ws.Range("E5").Value = 184361
$ws.Range("G5").Value = 152
$ws.Range("H5").Value = 22165

# Row 67: updated totals for this country
$ws.Range("B67").Value = 5915
$ws.Range("C67").Value = 336
$ws.Range("D67").Value = 609
$ws.Range("E67").Value = 5066
$ws.Range("G67").Value = 10
$ws.Range("H67").Value = 240

# Venezuela is newly inserted into the ranking at row 110; the countries
# that were there before (Albania, Guinea Ecuatorial) shift down one row
# each, keeping their own figures unchanged. Niger (row 113) is untouched.
$ws.Range("A110").Value = "Venezuela"
$ws.Range("B110").Value = 1010
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 262
$ws.Range("E110").Value = 738
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 10

$ws.Range("A111").Value = "Albania"
$ws.Range("B111").Value = 989
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 783
$ws.Range("E111").Value = 175
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 31

$ws.Range("A112").Value = "Guinea Ecuatorial"
$ws.Range("B112").Value = 960
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 165
$ws.Range("E112").Value = 784
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 11
